$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text so purely-numeric-looking strings
# (e.g. "305.71") are stored as text, matching the source data which
# uses inline/string cells throughout column D. Style is reset to
# "Normal" afterwards so no stray cell-style index is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '41.926.16'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '2.273.52'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '305.71'
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('D6').Value = '93.21'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('D10').Value = '32.77'
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('D11').Value = '0.0799'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('E12').Value = '  -1.81%  '
$ws.Range('D13').Value = '6.69'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').Value = '2.628.47'
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Value = '2.269.04'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('E17').Value = '  +3.39%  '
$ws.Range('D18').Value = '41.864.17'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = '12.74'
$ws.Range('E19').Value = '  +4.60%  '
$ws.Range('D20').Value = '0.0₃0915'
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '68.09'
$ws.Range('E22').Value = '  +1.29%  '
$ws.Range('D23').Value = '244.08'
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D24').Value = '2.60'
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('D25').Value = '1.95'
$ws.Range('E25').Value = '  +2.03%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '24.01'
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').Value = '9.70'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').Value = '35.33'
$ws.Range('E30').Value = '  +3.62%  '
$ws.Range('D31').Value = '159.26'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').Value = '5.38'
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '0.0744'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('D35').Value = '3.03'
$ws.Range('E35').Value = '  -0.93%  '
$ws.Range('D36').Value = '17.26'
$ws.Range('E36').Value = '  +3.58%  '
$ws.Range('E37').Value = '  -1.04%  '
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('E39').Value = '  +0.70%  '
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').Value = '19.79'
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('D43').Value = '2.012.23'
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('D44').Value = '2.25'
$ws.Range('E44').Value = '  +10.71%  '
$ws.Range('E45').Value = '  +1.15%  '
$ws.Range('D46').Value = '10.29'
$ws.Range('E46').Value = '  +1.74%  '
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('D48').Value = '53.54'
$ws.Range('E48').Value = '  +3.10%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '1.51'
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').Value = '72.61'
$ws.Range('E50').Value = '  +2.94%  '
$ws.Range('E51').Value = '  +0.29%  '

$ws.Range("D2:D51").Style = "Normal"
